$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:r/>" +
           "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
           "<w:r><w:t>: Join Aldo's Adventure across 4 different scenarios in this captivating 5x5 grid slot machine. Play for free and discover unique bonuses.</w:t></w:r>" +
           "</w:p>" +
           "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'></w:p>"

$insertPoint.InsertXML($metaXml) | Out-Null

# The snippet above is split across two <w:p> elements so that Word
# actually creates a paragraph break after our new paragraph instead of
# merging its runs into the following paragraph. That leaves a spare
# empty paragraph behind the new one which we now remove.
$spurious = $d.Paragraphs(3)
$spurious.Range.Delete()

# ------------------------------------------------------------------
# 2) Remove the duplicated title paragraph ("Play Aldo's Journey for
#    Free - Exciting Gameplay and Unique Bonuses") that used to sit
#    before the closing "Join Aldo's Adventure..." paragraph. Search
#    for it starting right after the real (first) title paragraph so
#    the original one is never touched.
# ------------------------------------------------------------------
$dupTitleRange = $d.Range($titlePara.Range.End, $d.Content.End)
$dupTitleRange.Find.Execute("Play Aldo's Journey for Free - Exciting Gameplay and Unique Bonuses", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dupTitleRange.Paragraphs(1).Range.Delete()

# ------------------------------------------------------------------
# 3) Turn the closing italic paragraph's text from the old "Join
#    Aldo's Adventure..." blurb into the new image-generation prompt.
#    This is the very last paragraph in the document, so address it
#    positionally to avoid matching the similar text that now also
#    appears in the new "Meta description" paragraph near the top.
# ------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$closingRange = $d.Paragraphs($lastParaIndex).Range.Duplicate
$closingRange.Find.Execute("Join Aldo's Adventure across 4 different scenarios in this captivating 5x5 grid slot machine. Play for free and discover unique bonuses.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closingRange.Text = "Prompt: Create a feature image for Aldo's Journey that captures the game's adventurous and fun theme. The image should be in a cartoon style and feature a happy Maya warrior with glasses."
